$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition list) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 380
$ws1.Range("F5").Value = 426
$ws1.Range("F6").Value = 269
$ws1.Range("F7").Value = 2427
$ws1.Range("F8").Value = 420
$ws1.Range("F9").Value = 6362
$ws1.Range("F10").Value = 167
$ws1.Range("F11").Value = 411
$ws1.Range("F12").Value = 23

# --- Sheet "全部类型" (all types, combined list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 380
$ws4.Range("F5").Value = 426
$ws4.Range("F6").Value = 269
$ws4.Range("F9").Value = 2427
$ws4.Range("F10").Value = 420
$ws4.Range("F11").Value = 6362
$ws4.Range("F12").Value = 167
$ws4.Range("F13").Value = 411
$ws4.Range("F15").Value = 23
